$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the 3rd data row (formerly row 3) down to row 4, leaving row 3 blank
# (support for n_blank_allowed param to TabularIo.load).
$ws.Range("A3:D3").Cut($ws.Range("A4:D4"))
$ws.Range("A3:D3").Clear()

# The hyperlink that was anchored on A3 needs to follow its data down to A4.
# Cut/Paste does not re-anchor hyperlinks automatically, so remove the old
# one and re-add it at the new location, restoring the original Hyperlink
# cell style afterwards.
$savedStyle = $ws.Range("A4").Style
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$3') {
        $hl.Delete()
    }
}
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:foo@bar.com") | Out-Null
$ws.Range("A4").Style = $savedStyle

# Match the new selection left behind in the sheet.
$ws.Range("E3").Select()
